# "isEnd to pos integer" - change the True/False "isEnd" label on the
# tree-draw node ovals (slide 2) to the integer form (1/0) that matches
# the "pos" row already used above it.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# Shape 1 ("Oval 10", inside "Group 22") only needs its 3rd paragraph's
# text swapped (False -> 0); its run keeps its existing language/run
# properties and the paragraph keeps its endParaRPr, same as the diff.
$group22 = $s.Shapes.Item(2).GroupItems
$oval10 = $group22.Item(1)
$oval10.TextFrame.TextRange.Paragraphs(3).Text = "0"

# The remaining six ovals are fully retyped (all three lines) so the
# trailing paragraph loses its now-redundant endParaRPr, matching the
# rest of the diff's hunks.
function Set-NodeText($shape, $line1, $line2, $line3) {
    $tf = $shape.TextFrame
    $tr = $tf.TextRange
    $tr.Delete()
    $tf.TextRange.Text = $line1 + [char]13 + $line2 + [char]13 + $line3
}

Set-NodeText $group22.Item(2) "10,6" "0" "1"   # Oval 16: True -> 1
Set-NodeText $group22.Item(3) "9,3" "1" "0"    # Oval 17: False -> 0
Set-NodeText $group22.Item(4) "14,1" "1" "1"   # Oval 18: True -> 1

$group1 = $s.Shapes.Item(3).GroupItems
Set-NodeText $group1.Item(1) "10,6" "0" "1"    # Oval 23: True -> 1
Set-NodeText $group1.Item(2) "9,3" "1" "0"     # Oval 24: False -> 0
Set-NodeText $group1.Item(3) "14,1" "1" "1"    # Oval 25: True -> 1
